$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.786.22'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.635.20'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = "'215.21"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').Value = "'0.506"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  +0.40%  '
$ws.Range('E9').Value = '  -0.69%  '
$ws.Range('D10').Value = "'19.62"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.88%  '
$ws.Range('D11').Value = "'0.0788"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.635.58'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.859.84'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('D16').Value = '0.0₃0768'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').Value = "'62.89"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').Value = '25.790.17'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = "'4.45"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('D21').Value = "'194.46"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('D23').Value = "'6.29"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.01%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E25').Value = '  +4.07%  '
$ws.Range('E26').Value = '  +3.07%  '
$ws.Range('E27').Value = '  -0.49%  '
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('D31').Value = "'0.0494"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('E33').Value = '  +0.23%  '
$ws.Range('E34').Value = '  +1.06%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = "'0.904"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').Value = '1.130.35'
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('E38').Value = '  -1.85%  '
$ws.Range('E39').Value = '  -1.39%  '
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('E42').Value = '  +2.37%  '
$ws.Range('D43').Value = "'100.28"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.06%  '
$ws.Range('D44').Value = "'0.806"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.70%  '
$ws.Range('D45').Value = '1.769.34'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').Value = '0.0₆0109'
$ws.Range('E46').Value = '  -3.19%  '
$ws.Range('D47').Value = "'55.20"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('D49').Value = "'0.0504"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').Value = "'7.54"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.20%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = "'2.32"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.43%  '
